# Updated cryptos list refresh: new Price (D) and Volume(1h) (E) values for
# each coin row. Price cells are text (e.g. "26.406.59", "0.568") so we
# prefix with a leading apostrophe to force Excel to keep them as text
# instead of auto-coercing to a number, then ClearFormats() to drop the
# "quote prefix" cell style that the apostrophe entry adds - this keeps the
# cell's style identical to the untouched cells (no explicit style index).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.406.59"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "'1.617.66"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'212.85"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").Value = "'19.18"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "'1.845.34"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("D13").Value = "'1.617.91"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("D14").Value = "'4.02"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "'63.93"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "'236.96"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +10.02%  "
$ws.Range("D18").Value = "'26.420.92"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").Value = "'7.84"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.25%  "
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").Value = "'4.30"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("E24").Value = "  +3.13%  "
$ws.Range("D25").Value = "'147.09"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").Value = "'15.50"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "'1.528.93"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +7.11%  "
$ws.Range("D33").Value = "'3.24"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("E35").Value = "  +4.61%  "
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("D37").Value = "'0.568"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").Value = "'0.832"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("E40").Value = "  +3.03%  "
$ws.Range("E42").Value = "  +1.76%  "
$ws.Range("D43").Value = "'1.756.68"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.48%  "
$ws.Range("D44").Value = "'0.765"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("D45").Value = "'61.65"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.49%  "
$ws.Range("D46").Value = "'0.906"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("D47").Value = "'90.80"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.76%  "
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("D50").Value = "'0.0962"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("D51").Value = "'7.52"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.63%  "
